$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Костырев Александр (row 12): ДЗ_1 (column C) was 0, now submitted -> 5,
# cleared to white fill but flagged in red text.
$ws.Range("C12").Value = 5
$ws.Range("C12").Font.Color = 255
$ws.Range("C12").Interior.ThemeColor = 2

# Гоман Антон (row 7): ДЗ_2 (column D) was 0, now submitted -> 5,
# highlight cleared back to white (done).
$ws.Range("D7").Value = 5
$ws.Range("D7").Interior.ThemeColor = 2

# Restore the working selection/scroll position used while editing.
$ws.Range("D7").Select()

# Page setup for printing.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
